$d = $word.ActiveDocument

# Locate the paragraph ending in "... summarized in the table below:"
$rng = $d.Content
$rng.Find.Execute("Some of the analyst", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$srcPara = $rng.Paragraphs(1)
$srcRange = $srcPara.Range
$srcRange.Collapse(0)
$insertPos = $srcRange.End

# Insert a new empty paragraph right after it.
$srcRange.InsertParagraphAfter() | Out-Null

# Grab the freshly inserted (still empty) paragraph via its known position.
$newRange = $d.Range($insertPos, $insertPos)
$newPara = $newRange.Paragraphs(1)
$newPara.Style = "SourceCode"

# First run: "pander" styled as FunctionTok.
$codeRange = $newPara.Range
$codeRange.Collapse(0)
$codeRange.InsertAfter("pander")
$codeRange.Style = "FunctionTok"

# Second run: "(table_forecasts)" styled as NormalTok.
$codeRange.Collapse(0)
$codeRange.InsertAfter("(table_forecasts)")
$codeRange.Style = "NormalTok"
